{"js": "// Reorder \"havainnointijaksot vuonna Perseuksen t\u00e4hdist\u00f6 2022: ...\" to\n// \"Perseuksen t\u00e4hdist\u00f6 havainnointijaksot vuonna 2022: ...\" everywhere it\n// occurs in the document body (there are 4 identical occurrences).\nconst OLD_TEXT = \"havainnointijaksot vuonna Perseuksen t\u00e4hdist\u00f6 2022: 16.-25.1., 7.-16.11., 6.-15.12\";\nconst NEW_TEXT = \"Perseuksen t\u00e4hdist\u00f6 havainnointijaksot vuonna 2022: 16.-25.1., 7.-16.11., 6.-15.12\";\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nfor (const paragraph of paragraphs.items) {\n  if (paragraph.text === OLD_TEXT) {\n    paragraph.insertText(NEW_TEXT, Word.InsertLocation.replace);\n  }\n}\n\nawait context.sync();\n", "ps1": "# Reorder \"havainnointijaksot vuonna Perseuksen t\u00e4hdist\u00f6 2022: ...\" to\n# \"Perseuksen t\u00e4hdist\u00f6 havainnointijaksot vuonna 2022: ...\" everywhere it\n# occurs in the document body (there are 4 identical occurrences).\n$d = $word.ActiveDocument\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n\n$find.Text = \"havainnointijaksot vuonna Perseuksen t\u00e4hdist\u00f6 2022: 16.-25.1., 7.-16.11., 6.-15.12\"\n$find.Replacement.Text = \"Perseuksen t\u00e4hdist\u00f6 havainnointijaksot vuonna 2022: 16.-25.1., 7.-16.11., 6.-15.12\"\n$find.Forward = $true\n$find.Wrap = 1\n$find.Format = $false\n$find.MatchCase = $true\n$find.MatchWholeWord = $false\n$find.MatchWildcards = $false\n\n$find.Execute([ref]$find.Text, [ref]$find.MatchCase, [ref]$find.MatchWholeWord, [ref]$find.MatchWildcards, $null, $null, [ref]$find.Forward, $null, $null, [ref]$find.Replacement.Text, 2)\n"}
